# Delete the 2007年 data row (row 2). Excel shifts the remaining rows
# (2010年, 2012年, 2015年, 2017年) up by one, turning the former A1:R6
# used range into A1:R5 — matching the target diff exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(2).Delete()
